$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# --- Row 1: add a header-style data row identical to row 2's pattern ---
# (row 1 previously existed only as an empty, custom-formatted row)
$ws.Rows.Item(1).ClearFormats()
$ws.Range("A2:D2").Copy()
$ws.Range("A1:D1").PasteSpecial($xlPasteValues)

# --- Row 5: column B changes from "  102" (leading spaces) to plain "102" ---
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial($xlPasteValues)

# --- Helper cell used to stash numeric-looking text ("123") in the shared
#     string table as a genuine Text value without ever touching
#     NumberFormat/style (keeps styles.xml untouched). A formula that
#     evaluates to a text result is pasted-as-value, which yields plain
#     text without flipping the "quote prefix" style flag. ---
$helper = $ws.Cells.Item(50, 10)
$helper.Formula = '="123"'

# --- Rows 6-9: Services test case rows ---
foreach ($r in 6..9) {
    $helper.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteValues)  # A: 123
    $ws.Cells.Item($r, 2).Value2 = "test"                # B: test
    $ws.Cells.Item($r, 3).Value2 = "G4234"               # C: G4234

    # D column: must end up as a present-but-empty cell (no value/type),
    # matching the diff's self-closed <c r="Dn"/>. Clearing contents drops
    # the cell entirely in this engine, so nudge a format property back to
    # its own current value to force the (now blank) cell to persist.
    $d = $ws.Cells.Item($r, 4)
    $d.ClearContents()
    $d.Font.Bold = $d.Font.Bold
}

# --- Row 10: Utils test case row (B/C swapped vs rows 6-9, plus a D value) ---
$helper.Copy()
$ws.Cells.Item(10, 1).PasteSpecial($xlPasteValues)  # A: 123
$ws.Cells.Item(10, 2).Value2 = "G4234"              # B: G4234
$ws.Cells.Item(10, 3).Value2 = "test"               # C: test
$ws.Cells.Item(10, 4).Value2 = "2019-09-12T12:01:20.457Z"  # D: timestamp

# --- Clean up the helper cell so it leaves no trace in the saved sheet ---
$helper.Clear()
